$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 header: rename variant labels, clear the 3rd (unused) variant ---
$ws.Range("E2").Value = "Wireframe Florian"
$ws.Range("G2").Value = "Wireframe Noel"
$ws.Range("I2").Value = $null

# --- Row 3: clear the now-unused "Teilnutzwert" label above column I ---
$ws.Range("I3").Value = $null

# --- Align the highlight formatting in column D (rows 7 and 9 get the
#     same fill/border as the other input cells, row 6, in that column) ---
$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Weighting (column C) ---
$ws.Range("C4").Value = 0.25
$ws.Range("C5").Value = 0.15
$ws.Range("C6").Value = 0.35
$ws.Range("C7").Value = 0.05
$ws.Range("C8").Value = 0.15
$ws.Range("C9").Value = 0.05

# --- Ratings for "Wireframe Florian" (column D) ---
$ws.Range("D4").Value = 7
$ws.Range("D5").Value = 5
$ws.Range("D6").Value = 8
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 4
$ws.Range("D9").Value = 8

# --- Ratings for "Wireframe Noel" (column F) ---
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = 7
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = 3
$ws.Range("F8").Value = 9
$ws.Range("F9").Value = 7

# --- Fix the per-row formulas in E/G so they reference their own row's weight ---
$ws.Range("E4").Formula = "=C`$4*D4"
$ws.Range("E5").Formula = "=C`$5*D5"
$ws.Range("E6").Formula = "=C`$6*D6"
$ws.Range("E7").Formula = "=C`$7*D7"
$ws.Range("E8").Formula = "=C`$8*D8"
$ws.Range("E9").Formula = "=C`$9*D9"

$ws.Range("G4").Formula = "=C`$4*F4"
$ws.Range("G5").Formula = "=C`$5*F5"
$ws.Range("G6").Formula = "=C`$6*F6"
$ws.Range("G7").Formula = "=C`$7*F7"
$ws.Range("G8").Formula = "=C`$8*F8"
$ws.Range("G9").Formula = "=C`$9*F9"

# --- Selection moves to K8 ---
$ws.Range("K8").Select()

$wb.Save()
